$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-16 Monday" "2026-02-17 Tuesday"
Replace-Text "18+29=47" "26+60=86"
Replace-Text "74-55=19" "69-2=67"
Replace-Text "97-71=26" "65-9=56"
Replace-Text "77-36=41" "1+37=38"
Replace-Text "62+29=91" "78-35=43"
Replace-Text "46+44=90" "77-21=56"
Replace-Text "3+52=55" "17-15=2"
Replace-Text "90-8=82" "5+78=83"
Replace-Text "28+25=53" "3+10=13"
Replace-Text "68-34=34" "56-51=5"
Replace-Text "46+6=52" "12+54=66"
Replace-Text "84-47=37" "51+48=99"
Replace-Text "14-1=13" "1+95=96"
Replace-Text "67-60=7" "44-7=37"
Replace-Text "38+42=80" "6+48=54"
Replace-Text "27+70=97" "11+32=43"
Replace-Text "87-13=74" "48+41=89"
Replace-Text "15+28=43" "1+20=21"
Replace-Text "62-53=9" "87+5=92"
Replace-Text "87-60=27" "66-62=4"
Replace-Text "79-33=46" "0+37=37"
Replace-Text "86-79=7" "32+37=69"
Replace-Text "65+16=81" "4+75=79"
Replace-Text "22+35=57" "75+10=85"
Replace-Text "32+47=79" "64+27=91"
Replace-Text "84-46=38" "0+57=57"
Replace-Text "42+54=96" "48+30=78"
Replace-Text "13+68=81" "82-82=0"
Replace-Text "81-45=36" "91-65=26"
Replace-Text "9+24=33" "16+68=84"
Replace-Text "73-0=73" "0-0=0"
Replace-Text "6+15=21" "21-17=4"
Replace-Text "0+60=60" "49-6=43"
Replace-Text "15+38=53" "68-6=62"
Replace-Text "72-57=15" "99-40=59"
Replace-Text "75-31=44" "71+4=75"
Replace-Text "43+53=96" "1+37=38"
Replace-Text "39-8=31" "78-59=19"
Replace-Text "72-8=64" "1+72=73"
Replace-Text "29-18=11" "69-5=64"
Replace-Text "17+60=77" "64-53=11"
Replace-Text "23+34=57" "37-0=37"
Replace-Text "64+13=77" "59+10=69"
Replace-Text "3+22=25" "45+17=62"
Replace-Text "90-47=43" "3+39=42"
Replace-Text "96-14=82" "46+0=46"
Replace-Text "66-48=18" "31+31=62"
Replace-Text "29-16=13" "62+34=96"
Replace-Text "89-17=72" "75+18=93"
Replace-Text "82-31=51" "40+7=47"
Replace-Text "43+20=63" "55+31=86"
Replace-Text "39+48=87" "42-29=13"
Replace-Text "96-71=25" "10+48=58"
Replace-Text "72-40=32" "8+56=64"
Replace-Text "78-47=31" "61-31=30"
Replace-Text "41+36=77" "14-3=11"
Replace-Text "70-11=59" "3+85=88"
Replace-Text "45+5=50" "49+17=66"
Replace-Text "41+33=74" "65-37=28"
Replace-Text "19+61=80" "95-27=68"
Replace-Text "45+29=74" "31-29=2"
Replace-Text "3+87=90" "21+70=91"
Replace-Text "72+7=79" "20+23=43"
Replace-Text "15+35=50" "22+58=80"
Replace-Text "23-6=17" "29-25=4"
Replace-Text "89-37=52" "17+61=78"
Replace-Text "94-84=10" "63-45=18"
Replace-Text "82+7=89" "32-9=23"
Replace-Text "24+62=86" "88-56=32"
Replace-Text "15+12=27" "86-50=36"
Replace-Text "11+22=33" "23-7=16"
Replace-Text "30+67=97" "88-8=80"
Replace-Text "90-54=36" "50-28=22"
Replace-Text "27-19=8" "60+32=92"
Replace-Text "99-63=36" "67-62=5"
Replace-Text "75-68=7" "65-42=23"
Replace-Text "51+13=64" "66+27=93"
Replace-Text "70-57=13" "90+0=90"
Replace-Text "8+14=22" "42-28=14"
Replace-Text "52-30=22" "83-13=70"
Replace-Text "40-11=29" "92-12=80"
Replace-Text "15+1=16" "53+1=54"
Replace-Text "42-30=12" "8+63=71"
Replace-Text "42-16=26" "1+28=29"
Replace-Text "53+28=81" "96-12=84"
Replace-Text "26-15=11" "65-29=36"
Replace-Text "8+59=67" "21+42=63"
Replace-Text "19+74=93" "90-22=68"
Replace-Text "20+40=60" "80-66=14"
Replace-Text "38+39=77" "81-38=43"
Replace-Text "74-38=36" "83-20=63"
Replace-Text "93-49=44" "54+4=58"
Replace-Text "38+23=61" "69-44=25"
Replace-Text "41+31=72" "61-10=51"
Replace-Text "54-48=6" "45-25=20"
Replace-Text "35+32=67" "83+4=87"
Replace-Text "18-8=10" "42-5=37"
Replace-Text "37-2=35" "45-0=45"
Replace-Text "46-29=17" "18+17=35"
Replace-Text "5+13=18" "19+68=87"
